{"js": "// Swap the \"Nhi\u1ec7m v\u1ee5\" (task/description) text between the two data rows\n// of the team-evaluation table (STT | MSSV | H\u1ecd v\u00e0 T\u00ean | Nhi\u1ec7m v\u1ee5 | \u0110\u00e1nh gi\u00e1):\n//   Row 1 (L\u00ea Ho\u00e0ng Huy):        \"Th\u1ef1c hi\u1ec7n nh\u00f3m ch\u1ee9c n\u0103ng c\u01a1 b\u1ea3n c\u1ee7a user, nh\u00f3m ch\u1ee9c n\u0103ng qu\u1ea3n l\u00ed ng\u01b0\u1eddi d\u00f9ng c\u1ee7a admin\"\n//   Row 2 (L\u00ea Quang \u0110\u0103ng Khoa):  \"Nh\u00f3m ch\u1ee9c n\u0103ng v\u1ec1 v\u00ed, nh\u00f3m ch\u1ee9c n\u0103ng qu\u1ea3n l\u00ed v\u00ed c\u1ee7a admin\"\n// After the edit these two values are swapped between the rows.\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\n// Locate the table that has the expected header row: STT / MSSV / H\u1ecd v\u00e0 T\u00ean / Nhi\u1ec7m v\u1ee5 / \u0110\u00e1nh gi\u00e1\nfor (let i = 0; i < tables.items.length; i++) {\n  tables.items[i].load(\"values\");\n}\nawait context.sync();\n\nlet target = null;\nfor (let i = 0; i < tables.items.length; i++) {\n  const t = tables.items[i];\n  const v = t.values;\n  if (\n    v.length > 0 &&\n    v[0].length >= 4 &&\n    v[0][0].trim() === \"STT\" &&\n    v[0][3].trim() === \"Nhi\u1ec7m v\u1ee5\"\n  ) {\n    target = t;\n    break;\n  }\n}\n\nif (!target) {\n  throw new Error(\"Could not locate the team-evaluation table.\");\n}\n\n// Column index 3 is \"Nhi\u1ec7m v\u1ee5\"; row index 1 and 2 are the two data rows.\nconst cellA = target.getCell(1, 3);\nconst cellB = target.getCell(2, 3);\n\n// Use the cell's (single) paragraph range so the existing run/paragraph\n// formatting (and identity) is preserved, instead of rebuilding the cell\n// body from scratch the way `cell.value = ...` would.\nconst paraA = cellA.body.paragraphs.getFirst();\nconst paraB = cellB.body.paragraphs.getFirst();\nconst rangeA = paraA.getRange();\nconst rangeB = paraB.getRange();\nrangeA.load(\"text\");\nrangeB.load(\"text\");\nawait context.sync();\n\nconst textA = rangeA.text;\nconst textB = rangeB.text;\n\nrangeA.insertText(textB, Word.InsertLocation.replace);\nrangeB.insertText(textA, Word.InsertLocation.replace);\nawait context.sync();\n", "ps1": "# Swap the \"Nhi\u1ec7m v\u1ee5\" (task/description) text between the two data rows\n# of the team-evaluation table (STT | MSSV | H\u1ecd v\u00e0 T\u00ean | Nhi\u1ec7m v\u1ee5 | \u0110\u00e1nh gi\u00e1):\n#   Row 1 (L\u00ea Ho\u00e0ng Huy):        \"Th\u1ef1c hi\u1ec7n nh\u00f3m ch\u1ee9c n\u0103ng c\u01a1 b\u1ea3n c\u1ee7a user, nh\u00f3m ch\u1ee9c n\u0103ng qu\u1ea3n l\u00ed ng\u01b0\u1eddi d\u00f9ng c\u1ee7a admin\"\n#   Row 2 (L\u00ea Quang \u0110\u0103ng Khoa):  \"Nh\u00f3m ch\u1ee9c n\u0103ng v\u1ec1 v\u00ed, nh\u00f3m ch\u1ee9c n\u0103ng qu\u1ea3n l\u00ed v\u00ed c\u1ee7a admin\"\n# After the edit these two values are swapped between the rows.\n\n$d = $word.ActiveDocument\n\n# Locate the table whose header row is STT / MSSV / H\u1ecd v\u00e0 T\u00ean / Nhi\u1ec7m v\u1ee5 / \u0110\u00e1nh gi\u00e1\n$target = $null\nfor ($i = 1; $i -le $d.Tables.Count; $i++) {\n    $t = $d.Tables.Item($i)\n    if ($t.Columns.Count -eq 5 -and $t.Rows.Count -ge 3) {\n        $header1 = $t.Cell(1, 1).Range.Text.TrimEnd([char]7, [char]13).Trim()\n        if ($header1 -eq \"STT\") {\n            $target = $t\n            break\n        }\n    }\n}\n\nif ($target -eq $null) {\n    throw \"Could not locate the team-evaluation table.\"\n}\n\n# Column 4 is \"Nhi\u1ec7m v\u1ee5\"; data rows are row 2 (L\u00ea Ho\u00e0ng Huy) and row 3 (L\u00ea Quang \u0110\u0103ng Khoa).\n$cellA = $target.Cell(2, 4)\n$cellB = $target.Cell(3, 4)\n\n$textA = $cellA.Range.Text.TrimEnd([char]7, [char]13)\n$textB = $cellB.Range.Text.TrimEnd([char]7, [char]13)\n\n$cellA.Range.Text = $textB\n$cellB.Range.Text = $textA\n"}
